$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W1").Value = 0.99781116221238064
$ws.Range("A2").Value = 0.7877880531143715
$ws.Range("BJ2").Value = 0.70104607801305674
$ws.Range("N3").Value = 0.87800688409249383
$ws.Range("BL3").Value = 0.940447228446322
$ws.Range("AX4").Value = 0.93780537555193122
$ws.Range("AZ4").Value = 0.99732103307627451
$ws.Range("BI4").Value = 0.90734391903226164
$ws.Range("F5").Value = 0.59675404024587397
$ws.Range("X5").Value = 0.96542120821948796
$ws.Range("AW6").Value = 0.77699721114876585
$ws.Range("AZ6").Value = 0.7645138071840617
$ws.Range("BL6").Value = 0.60966340164738586
$ws.Range("BO6").Value = 0.67845507932397364
$ws.Range("AO7").Value = 0.87769454161491911
$ws.Range("BJ7").Value = 0.76877567786286594
$ws.Range("J8").Value = 0.9397934689347176
$ws.Range("AG8").Value = 0.97026755900740591
$ws.Range("BM8").Value = 0.96789251101386642
$ws.Range("U10").Value = 0.99715683339200245
$ws.Range("O12").Value = 0.87010783537513081
$ws.Range("AM13").Value = 0.57159546214431545
$ws.Range("BG15").Value = 0.97181360312405274
$ws.Range("BH15").Value = 0.85799011891073063
$ws.Range("R16").Value = 0.93631183073855839
$ws.Range("BP17").Value = 0.91444686023393573
$ws.Range("Q18").Value = 0.75005020973634329
$ws.Range("I19").Value = 0.63634018608146015
$ws.Range("R19").Value = 0.85578381681088422
$ws.Range("U19").Value = 0.98106667409995363
$ws.Range("AI19").Value = 0.8405669711950372
$ws.Range("AU19").Value = 0.91590139726914666
$ws.Range("O20").Value = 0.93050270579383154
$ws.Range("P21").Value = 0.97882800243759727
$ws.Range("AH22").Value = 0.86474921170237828
$ws.Range("AJ22").Value = 0.71597655374859337
$ws.Range("BN22").Value = 0.72277250238081781
$ws.Range("L24").Value = 0.83366064252006666
$ws.Range("N24").Value = 0.97540844986313702
$ws.Range("W24").Value = 0.82560812331306688
$ws.Range("BL24").Value = 0.92404261932244414
$ws.Range("I25").Value = 0.59484434296133704
$ws.Range("AO25").Value = 0.89747057611836956
$ws.Range("BJ25").Value = 0.52585238186590311
$ws.Range("A26").Value = 0.98515979264737519
$ws.Range("P27").Value = 0.65880072822212565
$ws.Range("AH28").Value = 0.99481877150133835
$ws.Range("BO28").Value = 0.95651937860125913
$ws.Range("E29").Value = 0.88642792516529201
$ws.Range("K29").Value = 0.66556994729007735
$ws.Range("S29").Value = 0.96031681791181711
$ws.Range("BC29").Value = 0.77732768700582944
$ws.Range("P30").Value = 0.61353792470639923
$ws.Range("R31").Value = 0.8942266382194739
$ws.Range("AG31").Value = 0.92619630632589955
$ws.Range("AV31").Value = 0.88796448973398423
$ws.Range("AB32").Value = 0.95820259933637164
$ws.Range("BM32").Value = 0.72453214924973863
$ws.Range("G33").Value = 0.70391552384748213
$ws.Range("AS33").Value = 0.8982997025361803
$ws.Range("AA34").Value = 0.63276983550089727
$ws.Range("AF34").Value = 0.66342008123163576
$ws.Range("AO34").Value = 0.87693052547094896
$ws.Range("E35").Value = 0.89571757501775062
$ws.Range("I35").Value = 0.64004253559290436
$ws.Range("BI35").Value = 0.96233832767474858
$ws.Range("M36").Value = 0.61231170061288187
$ws.Range("AU36").Value = 0.93875075198276692
$ws.Range("BB36").Value = 0.71803893272768182
$ws.Range("AF37").Value = 0.99664165029539231
$ws.Range("AH37").Value = 0.89909211852955051
$ws.Range("AJ38").Value = 0.91744148224898148
$ws.Range("AP38").Value = 0.56831277262229984
$ws.Range("BE38").Value = 0.74718812210175301
$ws.Range("F39").Value = 0.87509089949499841
$ws.Range("J39").Value = 0.73039632190214854
$ws.Range("U39").Value = 0.95580300968634613
$ws.Range("AQ39").Value = 0.82625452847251291
$ws.Range("AR39").Value = 0.8661026341286574
$ws.Range("T40").Value = 0.89487337847717963
$ws.Range("AO40").Value = 0.55718836207830835
$ws.Range("AU40").Value = 0.9600527073262386
$ws.Range("K41").Value = 0.78985167012875934
$ws.Range("F42").Value = 0.92947352317665155
$ws.Range("AA42").Value = 0.89567274977399935
$ws.Range("AQ42").Value = 0.62477911807995334
$ws.Range("BH42").Value = 0.79166481351433016
$ws.Range("BI42").Value = 0.89174577106820418
$ws.Range("AP44").Value = 0.56204263983103475
$ws.Range("AR45").Value = 0.88801308083227903
$ws.Range("I46").Value = 0.74662917039478849
$ws.Range("R46").Value = 0.67204202091737708
$ws.Range("AS46").Value = 0.93978419566790794
$ws.Range("Q47").Value = 0.71108914066891016
$ws.Range("AT47").Value = 0.87595484692819614
$ws.Range("Q48").Value = 0.90898548486373876
$ws.Range("BM48").Value = 0.80232981542810244
$ws.Range("Z49").Value = 0.8758337087531558
$ws.Range("C50").Value = 0.9002131588729454
$ws.Range("AH50").Value = 0.97748139363058228
$ws.Range("AO51").Value = 0.98942614165161147
$ws.Range("AX51").Value = 0.67231994242975279
$ws.Range("AY52").Value = 0.83880760566905166
$ws.Range("D53").Value = 0.96523102118340109
$ws.Range("BE53").Value = 0.93363514143697834
$ws.Range("B54").Value = 0.91258304955379232
$ws.Range("O54").Value = 0.81963162051629035
$ws.Range("BF54").Value = 0.91698595871254207
$ws.Range("J55").Value = 0.95141153409119439
$ws.Range("AQ55").Value = 0.94583125771917698
$ws.Range("BD55").Value = 0.94888353994154873
$ws.Range("BG55").Value = 0.88330098217737107
$ws.Range("B56").Value = 0.92299787241570586
$ws.Range("AD56").Value = 0.67924811649757078
$ws.Range("AV56").Value = 0.892589871526418
$ws.Range("AX56").Value = 0.73886230805607311
$ws.Range("BJ57").Value = 0.97436336197659812
$ws.Range("C58").Value = 0.9904229801791431
$ws.Range("Z58").Value = 0.56248793393201446
$ws.Range("AG58").Value = 0.57851609716498986
$ws.Range("AH58").Value = 0.85996094859420724
$ws.Range("AI59").Value = 0.98941146369698019
$ws.Range("P60").Value = 0.91596166605859719
$ws.Range("U62").Value = 0.87919654675867098
$ws.Range("X62").Value = 0.86773203420568357
$ws.Range("AJ63").Value = 0.89259085341851274
$ws.Range("BF63").Value = 0.662810734528961
$ws.Range("BJ63").Value = 0.92896361914256653
$ws.Range("AD65").Value = 0.81234943273940741
$ws.Range("AR65").Value = 0.84767544164758291
$ws.Range("BA65").Value = 0.75559275026216932
$ws.Range("AW66").Value = 0.89321907965537928
$ws.Range("BP66").Value = 0.65763662477109164
$ws.Range("AT67").Value = 0.71819691495260285
$ws.Range("I68").Value = 0.98789388353232088
$ws.Range("U68").Value = 0.68235171167596698
